$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Planificación
$ws.Range("E3").Value = "Documentacion de estimaciones del proyecto"
$ws.Range("H3").Value = "enero,2023-febrero,2023"
$ws.Range("I3").Value = "Si se realiza muy deprisa o de manera no adecuado todas las demas fases tendran errores"

# Row 4 - Consultoria
$ws.Range("E4").Value = "Diseño funcional"
$ws.Range("H4").Value = "febrero,2023-marzo,2023"
$ws.Range("I4").Value = "Realizar esta fase incorrectamente, provocaria que errores se produjeran en las siguientes fases obligando a echar atras el proyecto para tener que volverlo a implementar"

# Row 5 - Diseño
$ws.Range("E5").Value = "Diseño tecnico"
$ws.Range("H5").Value = "marzo,2023-julio,2023"
$ws.Range("I5").Value = "Errores a la hora de realizar el diseño tecnico lo cual provocaria un funcionamiento inesperado"

# Row 6 - Pruebas
$ws.Range("E6").Value = "Documentacion de pruebas realizadas"
$ws.Range("H6").Value = "julio,2023-septiembre,2023"
$ws.Range("I6").Value = "No descubrir errores los cuales obligaran a volver a anterior fases de implementacion cuando se descubran mas adelante"

# Row 7 - Formacion
$ws.Range("H7").Value = "septiembre.2023-noviembre.2023"
$ws.Range("I7").Value = "Una formacion ineficiente provocara en que los empleados puedan realizar un uso inadecuado del ERP"

# Row 8 - Arranque
$ws.Range("H8").Value = "noviembre.2023-diciembre.2023"
$ws.Range("I8").Value = "Un arranque mal realizado puedo llevar a que la empresa tenga que detener su funcionamiento temporalmente"

# Row 9 - Soporte
$ws.Range("H9").Value = "diciembre,2023-enero.2023"
$ws.Range("I9").Value = "Un mal soporte provocara que el sistema pueda tener fallos que afecten a la empresa"

# Row 10 - Mantenimiento
$ws.Range("H10").Value = "enero,2023-"
$ws.Range("I10").Value = "Que cambios o actualizaciones realizadas provoquen fallos en el ERP los cuales impidan su debido funcionamiento"

# Update active cell selection to match the saved state of the workbook
$ws.Range("G10").Select()
